$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-08-08"

# Update header label for column B (the "through" month column)
$ws.Range("B1").Value = "August 2022 (through August 08)"

# Helper to set a cell value by A1 reference
function Set-Cell($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Row 2 - Austin
Set-Cell "J2" 3

# Row 3 - Grand Crossing
Set-Cell "AP3" 2

# Row 4 - North Lawndale
Set-Cell "J4" 2
Set-Cell "AP4" 4

# Row 5 - Garfield Park
Set-Cell "R5" 2
Set-Cell "Z5" 1

# Row 6 - Englewood
Set-Cell "B6" 4

# Row 7 - Humboldt Park
Set-Cell "B7" 3
Set-Cell "J7" 2

# Row 9 - Chatham
Set-Cell "R9" 3

# Row 10 - Little Village
Set-Cell "J10" 3

# Row 13 - Roseland
Set-Cell "J13" 3
Set-Cell "R13" 1

# Row 23 - Washington Heights
Set-Cell "J23" 1

# Row 27 - New City
Set-Cell "R27" 1

# Row 32 - Calumet Heights
Set-Cell "J32" 1
Set-Cell "AH32" 2

# Row 38 - South Deering
Set-Cell "J38" 1

# Row 41 - Rogers Park
Set-Cell "B41" 1

# Row 46 - Kenwood
Set-Cell "J46" 2

# Row 63 - Brighton Park
Set-Cell "R63" 1

# Row 64 - Bucktown
Set-Cell "R64" 2

# Row 65 - Burnside
Set-Cell "J65" 1

# Row 91 - Pullman
Set-Cell "B91" 1
